# Update the column header label in cell B2 from
# "Average salary (not adjusted)" to "Average salary".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Average salary"
